$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Grupo_Experimental (column B) values for rows 2-6
$ws.Range("B2").Value = "Con SmartScore"
$ws.Range("B3").Value = "Sin SmartScore"
$ws.Range("B4").Value = "Con SmartScore"
$ws.Range("B5").Value = "Sin SmartScore"
$ws.Range("B6").Value = "Con SmartScore"

# Convert SmartScore text values in row 6 to true numeric values
$ws.Range("I6").Value = 0.54
$ws.Range("L6").Value = 0.52
$ws.Range("O6").Value = 0.449
$ws.Range("R6").Value = 0.622
$ws.Range("U6").Value = 0.616
$ws.Range("X6").Value = 0.615
$ws.Range("AA6").Value = 0.729
$ws.Range("AD6").Value = 0.702
$ws.Range("AG6").Value = 0.685
